$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: update title (D26) and link (E26)
$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"
$ws.Range("E26").Value = "https://blog.est.ai/2022/02/%ec%83%9d%ec%84%b1-%eb%aa%a8%eb%8d%b8%ec%9d%98-%ec%83%88%eb%a1%9c%ec%9a%b4-%ed%9d%90%eb%a6%84-%ed%99%95%ec%82%b0-%eb%aa%a8%eb%8d%b8diffusion-model%ec%97%90-%ea%b4%80%ed%95%98%ec%97%ac/"

# Row 37: update title (D37) and link (E37)
$ws.Range("D37").Value = "[Paper Review] Autoformer: Decomposition transformers with auto-correlation for long-term series forecasting"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1930&mod=document&pageid=1"
